$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing data rows down by one).
$ws.Rows("2:2").Insert()

# The freshly inserted row inherits formatting from the header row above it
# (bold font, border, centered alignment). Strip that so the new data row
# starts out unstyled, matching the rest of the data rows.
$ws.Range("A2:R2").ClearFormats()

# Column D carries a custom date number format; copy just that formatting
# from the row below (the old row 2, now row 3) onto the new D2 cell.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# Populate the new weekly entry (week of 2022-09-08).
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C2").Value = 'Los Lagos'
$ws.Range("D2").Value = 44812
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 100112012
$ws.Range("G2").Value = 'Espinaca'
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("N2").Value = '$/cuna 10 kilos'
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 1200
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 'Hortaliza'
